# Update "想去人数" (want-to-go count) figures in F column across all sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 447
$ws.Range("F5").Value = 361
$ws.Range("F6").Value = 551
$ws.Range("F7").Value = 46
$ws.Range("F10").Value = 382
$ws.Range("F12").Value = 668
$ws.Range("F13").Value = 754
$ws.Range("F14").Value = 1512
$ws.Range("F15").Value = 1512
$ws.Range("F16").Value = 886
$ws.Range("F18").Value = 1350
$ws.Range("F19").Value = 161
$ws.Range("F20").Value = 318
$ws.Range("F24").Value = 6598
$ws.Range("F25").Value = 4948
$ws.Range("F26").Value = 143
$ws.Range("F27").Value = 488
$ws.Range("F28").Value = 207
$ws.Range("F29").Value = 178
$ws.Range("F32").Value = 1282
$ws.Range("F33").Value = 193
$ws.Range("F35").Value = 611
$ws.Range("F43").Value = 97

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 18
$ws.Range("F18").Value = 239

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2455
$ws.Range("F4").Value = 195
$ws.Range("F5").Value = 56

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 447
$ws.Range("F6").Value = 195
$ws.Range("F7").Value = 56
$ws.Range("F8").Value = 361
$ws.Range("F9").Value = 551
$ws.Range("F10").Value = 46
$ws.Range("F14").Value = 382
$ws.Range("F16").Value = 668
$ws.Range("F17").Value = 754
$ws.Range("F18").Value = 1512
$ws.Range("F19").Value = 1512
$ws.Range("F20").Value = 886
$ws.Range("F22").Value = 1350
$ws.Range("F23").Value = 161
$ws.Range("F24").Value = 318
$ws.Range("F29").Value = 6598
$ws.Range("F30").Value = 4948
$ws.Range("F31").Value = 143
$ws.Range("F33").Value = 1282
$ws.Range("F34").Value = 193
$ws.Range("F38").Value = 611
$ws.Range("F47").Value = 97
$ws.Range("F49").Value = 239
